$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows 239-356: shift weekly data down by one week, columns D,I,J,K,L,M,O,P ---
$updates = @(
  @{Row=239; D=44917; I='Primera'; J=2000; K=800; L=900; M=850; O='Región Metropolitana'; P=850},
  @{Row=240; D=44917; I='Segunda'; J=1000; K=700; L=700; M=700; O='Región Metropolitana'; P=700},
  @{Row=241; D=44336; I='Primera'; J=2000; K=600; L=700; M=650; O='Región Metropolitana'; P=650},
  @{Row=242; D=44336; I='Segunda'; J=1000; K=500; L=500; M=500; O='Región Metropolitana'; P=500},
  @{Row=243; D=44775; I='Primera'; J=2000; K=900; L=1000; M=950; O='Región Metropolitana'; P=950},
  @{Row=244; D=44782; I='Primera'; J=2000; K=900; L=1000; M=950; O='Región Metropolitana'; P=950},
  @{Row=245; D=44782; I='Segunda'; J=1000; K=800; L=800; M=800; O='Región Metropolitana'; P=800},
  @{Row=246; D=44894; I='Primera'; J=1000; K=700; L=800; M=750; O='Región Metropolitana'; P=750},
  @{Row=247; D=44894; I='Segunda'; J=500; K=600; L=600; M=600; O='Región Metropolitana'; P=600},
  @{Row=248; D=44727; I='Primera'; J=2000; K=700; L=800; M=750; O='Región Metropolitana'; P=750},
  @{Row=249; D=44750; I='Primera'; J=1600; K=1000; L=1200; M=1075; O='Región Metropolitana'; P=1075},
  @{Row=250; D=44680; I='Primera'; J=2000; K=900; L=1000; M=950; O='Región Metropolitana'; P=950},
  @{Row=251; D=44680; I='Segunda'; J=1000; K=800; L=800; M=800; O='Región Metropolitana'; P=800},
  @{Row=252; D=44705; I='Primera'; J=2000; K=900; L=1000; M=950; O='Región Metropolitana'; P=950},
  @{Row=253; D=44232; I='Primera'; J=1000; K=700; L=800; M=750; O='Región Metropolitana'; P=750},
  @{Row=254; D=44232; I='Segunda'; J=500; K=600; L=600; M=600; O='Región Metropolitana'; P=600},
  @{Row=255; D=44448; I='Primera'; J=1000; K=700; L=800; M=750; O='Región Metropolitana'; P=750},
  @{Row=256; D=44448; I='Segunda'; J=500; K=600; L=600; M=600; O='Región Metropolitana'; P=600},
  @{Row=257; D=44657; I='Primera'; J=2000; K=1000; L=1200; M=1100; O='Región Metropolitana'; P=1100},
  @{Row=258; D=44657; I='Segunda'; J=1000; K=800; L=800; M=800; O='Región Metropolitana'; P=800},
  @{Row=259; D=44498; I='Primera'; J=2500; K=600; L=650; M=630; O='Región Metropolitana'; P=630},
  @{Row=260; D=44685; I='Primera'; J=1800; K=800; L=900; M=844; O='Región Metropolitana'; P=844},
  @{Row=261; D=44685; I='Segunda'; J=1700; K=600; L=650; M=626; O='Región Metropolitana'; P=626},
  @{Row=262; D=44908; I='Primera'; J=2000; K=1000; L=1100; M=1050; O='Región Metropolitana'; P=1050},
  @{Row=263; D=44908; I='Segunda'; J=1000; K=900; L=900; M=900; O='Región Metropolitana'; P=900},
  @{Row=264; D=44600; I='Primera'; J=1500; K=650; L=700; M=677; O='Región Metropolitana'; P=677},
  @{Row=265; D=44460; I='Primera'; J=1000; K=600; L=700; M=650; O='Región Metropolitana'; P=650},
  @{Row=266; D=44460; I='Segunda'; J=500; K=500; L=500; M=500; O='Región Metropolitana'; P=500},
  @{Row=267; D=44777; I='Primera'; J=2000; K=900; L=1000; M=950; O='Región Metropolitana'; P=950},
  @{Row=268; D=44777; I='Segunda'; J=1500; K=700; L=700; M=700; O='Región Metropolitana'; P=700},
  @{Row=269; D=44532; I='Primera'; J=650; K=550; L=600; M=573; O='Región Metropolitana'; P=573},
  @{Row=270; D=44428; I='Primera'; J=1000; K=600; L=700; M=650; O='Región Metropolitana'; P=650},
  @{Row=271; D=44428; I='Segunda'; J=500; K=500; L=500; M=500; O='Región Metropolitana'; P=500},
  @{Row=272; D=44785; I='Primera'; J=2000; K=900; L=1000; M=950; O='Región Metropolitana'; P=950},
  @{Row=273; D=44785; I='Segunda'; J=1000; K=800; L=800; M=800; O='Región Metropolitana'; P=800},
  @{Row=274; D=44194; I='Primera'; J=800; K=800; L=900; M=850; O='Región Metropolitana'; P=850},
  @{Row=275; D=44194; I='Segunda'; J=400; K=700; L=700; M=700; O='Región Metropolitana'; P=700},
  @{Row=276; D=44322; I='Primera'; J=2000; K=600; L=700; M=650; O='Región del Maule'; P=650},
  @{Row=277; D=44322; I='Segunda'; J=1000; K=500; L=500; M=500; O='Región del Maule'; P=500},
  @{Row=278; D=44799; I='Primera'; J=2000; K=900; L=1000; M=950; O='Región Metropolitana'; P=950},
  @{Row=279; D=44799; I='Segunda'; J=1000; K=700; L=700; M=700; O='Región Metropolitana'; P=700},
  @{Row=280; D=44441; I='Primera'; J=1000; K=700; L=800; M=750; O='Región Metropolitana'; P=750},
  @{Row=281; D=44441; I='Segunda'; J=500; K=600; L=600; M=600; O='Región Metropolitana'; P=600},
  @{Row=282; D=44831; I='Primera'; J=310; K=1200; L=1300; M=1252; O='Región Metropolitana'; P=1252},
  @{Row=283; D=44831; I='Segunda'; J=300; K=1000; L=1000; M=1000; O='Región Metropolitana'; P=1000},
  @{Row=284; D=44474; I='Primera'; J=1000; K=600; L=700; M=650; O='Región Metropolitana'; P=650},
  @{Row=285; D=44474; I='Segunda'; J=500; K=500; L=500; M=500; O='Región Metropolitana'; P=500},
  @{Row=286; D=44651; I='Primera'; J=2200; K=1200; L=1300; M=1245; O='Región Metropolitana'; P=1245},
  @{Row=287; D=44761; I='Primera'; J=1800; K=800; L=900; M=844; O='Región Metropolitana'; P=844},
  @{Row=288; D=44882; I='Primera'; J=2000; K=600; L=700; M=650; O='Región Metropolitana'; P=650},
  @{Row=289; D=44882; I='Segunda'; J=1000; K=500; L=500; M=500; O='Región Metropolitana'; P=500},
  @{Row=290; D=44663; I='Primera'; J=2000; K=900; L=1000; M=950; O='Región del Maule'; P=950},
  @{Row=291; D=44663; I='Segunda'; J=1000; K=800; L=800; M=800; O='Región del Maule'; P=800},
  @{Row=292; D=44740; I='Primera'; J=900; K=700; L=800; M=756; O='Región Metropolitana'; P=756},
  @{Row=293; D=44895; I='Primera'; J=1000; K=700; L=800; M=750; O='Región Metropolitana'; P=750},
  @{Row=294; D=44895; I='Segunda'; J=500; K=600; L=600; M=600; O='Región Metropolitana'; P=600},
  @{Row=295; D=44386; I='Primera'; J=1000; K=700; L=800; M=750; O='Región Metropolitana'; P=750},
  @{Row=296; D=44386; I='Segunda'; J=500; K=600; L=600; M=600; O='Región Metropolitana'; P=600},
  @{Row=297; D=44425; I='Primera'; J=2000; K=600; L=700; M=650; O='Región Metropolitana'; P=650},
  @{Row=298; D=44425; I='Segunda'; J=1000; K=500; L=500; M=500; O='Región Metropolitana'; P=500},
  @{Row=299; D=44656; I='Primera'; J=1600; K=900; L=1000; M=994; O='Región Metropolitana'; P=994},
  @{Row=300; D=44292; I='Primera'; J=1000; K=700; L=800; M=750; O='Región Metropolitana'; P=750},
  @{Row=301; D=44617; I='Primera'; J=3000; K=700; L=800; M=767; O='Región del Maule'; P=767},
  @{Row=302; D=44400; I='Primera'; J=2000; K=700; L=800; M=750; O='Región Metropolitana'; P=750},
  @{Row=303; D=44400; I='Segunda'; J=1000; K=600; L=600; M=600; O='Región Metropolitana'; P=600},
  @{Row=304; D=44714; I='Primera'; J=2500; K=900; L=1000; M=940; O='Región Metropolitana'; P=940},
  @{Row=305; D=44776; I='Primera'; J=1000; K=900; L=1000; M=950; O='Región Metropolitana'; P=950},
  @{Row=306; D=44776; I='Segunda'; J=800; K=700; L=700; M=700; O='Región Metropolitana'; P=700},
  @{Row=307; D=44491; I='Primera'; J=1000; K=800; L=900; M=850; O='Región Metropolitana'; P=850},
  @{Row=308; D=44491; I='Segunda'; J=500; K=700; L=700; M=700; O='Región Metropolitana'; P=700},
  @{Row=309; D=44216; I='Primera'; J=1000; K=600; L=700; M=650; O='Región del Maule'; P=650},
  @{Row=310; D=44216; I='Segunda'; J=500; K=500; L=500; M=500; O='Región del Maule'; P=500},
  @{Row=311; D=44264; I='Primera'; J=1000; K=800; L=900; M=850; O='Región Metropolitana'; P=850},
  @{Row=312; D=44264; I='Segunda'; J=500; K=700; L=700; M=700; O='Región Metropolitana'; P=700},
  @{Row=313; D=44376; I='Primera'; J=2000; K=500; L=600; M=550; O='Región Metropolitana'; P=550},
  @{Row=314; D=44376; I='Segunda'; J=1000; K=400; L=400; M=400; O='Región Metropolitana'; P=400},
  @{Row=315; D=44847; I='Primera'; J=1500; K=1000; L=1000; M=1000; O='Región Metropolitana'; P=1000},
  @{Row=316; D=44847; I='Segunda'; J=1000; K=800; L=800; M=800; O='Región Metropolitana'; P=800},
  @{Row=317; D=44299; I='Primera'; J=1000; K=600; L=700; M=650; O='Región del Maule'; P=650},
  @{Row=318; D=44299; I='Segunda'; J=500; K=500; L=500; M=500; O='Región del Maule'; P=500},
  @{Row=319; D=44756; I='Primera'; J=2000; K=1200; L=1300; M=1250; O='Región Metropolitana'; P=1250},
  @{Row=320; D=44756; I='Segunda'; J=1000; K=1000; L=1000; M=1000; O='Región Metropolitana'; P=1000},
  @{Row=321; D=44364; I='Primera'; J=2000; K=600; L=700; M=650; O='Región Metropolitana'; P=650},
  @{Row=322; D=44364; I='Segunda'; J=1000; K=500; L=500; M=500; O='Región Metropolitana'; P=500},
  @{Row=323; D=44818; I='Primera'; J=2200; K=900; L=1000; M=955; O='Región Metropolitana'; P=955},
  @{Row=324; D=44453; I='Primera'; J=2000; K=600; L=700; M=650; O='Región Metropolitana'; P=650},
  @{Row=325; D=44453; I='Segunda'; J=1000; K=500; L=500; M=500; O='Región Metropolitana'; P=500},
  @{Row=326; D=44832; I='Primera'; J=1900; K=800; L=1000; M=895; O='Región Metropolitana'; P=895},
  @{Row=327; D=44516; I='Primera'; J=1500; K=500; L=550; M=523; O='Región Metropolitana'; P=523},
  @{Row=328; D=44162; I='Primera'; J=800; K=600; L=700; M=650; O='Región del Maule'; P=650},
  @{Row=329; D=44162; I='Segunda'; J=400; K=500; L=500; M=500; O='Región del Maule'; P=500},
  @{Row=330; D=44687; I='Primera'; J=2000; K=900; L=1000; M=950; O='Región Metropolitana'; P=950},
  @{Row=331; D=44533; I='Primera'; J=1500; K=500; L=600; M=547; O='Región Metropolitana'; P=547},
  @{Row=332; D=44665; I='Primera'; J=2000; K=700; L=800; M=750; O='Región Metropolitana'; P=750},
  @{Row=333; D=44665; I='Segunda'; J=1000; K=600; L=600; M=600; O='Región Metropolitana'; P=600},
  @{Row=334; D=44477; I='Primera'; J=2000; K=600; L=700; M=650; O='Región Metropolitana'; P=650},
  @{Row=335; D=44477; I='Segunda'; J=1000; K=500; L=500; M=500; O='Región Metropolitana'; P=500},
  @{Row=336; D=44813; I='Primera'; J=2000; K=1200; L=1300; M=1250; O='Región Metropolitana'; P=1250},
  @{Row=337; D=44813; I='Segunda'; J=1000; K=900; L=900; M=900; O='Región Metropolitana'; P=900},
  @{Row=338; D=44545; I='Primera'; J=2200; K=500; L=550; M=527; O='Región Metropolitana'; P=527},
  @{Row=339; D=44819; I='Primera'; J=1600; K=900; L=1000; M=950; O='Región Metropolitana'; P=950},
  @{Row=340; D=44295; I='Primera'; J=1000; K=700; L=800; M=750; O='Región Metropolitana'; P=750},
  @{Row=341; D=44295; I='Segunda'; J=500; K=600; L=600; M=600; O='Región Metropolitana'; P=600},
  @{Row=342; D=44630; I='Segunda'; J=900; K=1000; L=1200; M=1089; O='Región Metropolitana'; P=1089},
  @{Row=343; D=44741; I='Primera'; J=2000; K=800; L=900; M=850; O='Región Metropolitana'; P=850},
  @{Row=344; D=44741; I='Segunda'; J=1000; K=700; L=700; M=700; O='Región Metropolitana'; P=700},
  @{Row=345; D=44350; I='Primera'; J=2000; K=600; L=700; M=650; O='Región Metropolitana'; P=650},
  @{Row=346; D=44350; I='Segunda'; J=1000; K=500; L=500; M=500; O='Región Metropolitana'; P=500},
  @{Row=347; D=44890; I='Primera'; J=1600; K=700; L=750; M=719; O='Región Metropolitana'; P=719},
  @{Row=348; D=44763; I='Primera'; J=1500; K=1200; L=1300; M=1233; O='Región Metropolitana'; P=1233},
  @{Row=349; D=44763; I='Segunda'; J=1000; K=1000; L=1000; M=1000; O='Región Metropolitana'; P=1000},
  @{Row=350; D=44565; I='Primera'; J=2000; K=700; L=800; M=750; O='Región Metropolitana'; P=750},
  @{Row=351; D=44565; I='Segunda'; J=1000; K=600; L=600; M=600; O='Región Metropolitana'; P=600},
  @{Row=352; D=44589; I='Primera'; J=800; K=1000; L=1100; M=1038; O='Región Metropolitana'; P=1038},
  @{Row=353; D=44449; I='Primera'; J=1000; K=700; L=800; M=750; O='Región Metropolitana'; P=750},
  @{Row=354; D=44449; I='Segunda'; J=500; K=600; L=600; M=600; O='Región Metropolitana'; P=600},
  @{Row=355; D=44357; I='Primera'; J=2000; K=600; L=700; M=650; O='Región Metropolitana'; P=650},
  @{Row=356; D=44357; I='Segunda'; J=1000; K=500; L=500; M=500; O='Región Metropolitana'; P=500},
)

foreach ($u in $updates) {
  $r = $u.Row
  $ws.Cells.Item($r, 4).Value = $u.D
  $ws.Cells.Item($r, 9).Value = $u.I
  $ws.Cells.Item($r, 10).Value = $u.J
  $ws.Cells.Item($r, 11).Value = $u.K
  $ws.Cells.Item($r, 12).Value = $u.L
  $ws.Cells.Item($r, 13).Value = $u.M
  $ws.Cells.Item($r, 15).Value = $u.O
  $ws.Cells.Item($r, 16).Value = $u.P
}

# --- Append new rows 357-358 (full rows, copied forward from old 355/356) ---
$newRows = @(
  @{Row=357; A=11; B='Vega Monumental Concepción'; C='Bíobío'; D=44736; E=8; F=100112008; G='Coliflor'; H='Sin especificar'; I='Primera'; J=2000; K=800; L=900; M=850; N='$/unidad'; O='Región Metropolitana'; P=850; Q=1; R='Hortaliza'},
  @{Row=358; A=11; B='Vega Monumental Concepción'; C='Bíobío'; D=44736; E=8; F=100112008; G='Coliflor'; H='Sin especificar'; I='Segunda'; J=1000; K=700; L=700; M=700; N='$/unidad'; O='Región Metropolitana'; P=700; Q=1; R='Hortaliza'},
)

foreach ($nr in $newRows) {
  $r = $nr.Row
  $ws.Cells.Item($r, 1).Value = $nr.A
  $ws.Cells.Item($r, 2).Value = $nr.B
  $ws.Cells.Item($r, 3).Value = $nr.C
  $ws.Cells.Item($r, 4).Value = $nr.D
  $ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item($r - 1, 4).NumberFormat
  $ws.Cells.Item($r, 5).Value = $nr.E
  $ws.Cells.Item($r, 6).Value = $nr.F
  $ws.Cells.Item($r, 7).Value = $nr.G
  $ws.Cells.Item($r, 8).Value = $nr.H
  $ws.Cells.Item($r, 9).Value = $nr.I
  $ws.Cells.Item($r, 10).Value = $nr.J
  $ws.Cells.Item($r, 11).Value = $nr.K
  $ws.Cells.Item($r, 12).Value = $nr.L
  $ws.Cells.Item($r, 13).Value = $nr.M
  $ws.Cells.Item($r, 14).Value = $nr.N
  $ws.Cells.Item($r, 15).Value = $nr.O
  $ws.Cells.Item($r, 16).Value = $nr.P
  $ws.Cells.Item($r, 17).Value = $nr.Q
  $ws.Cells.Item($r, 18).Value = $nr.R
}

"Done."